# cnn_inference_trained.xlsx — move the highlighted "filter window" on the
# Inputs sheet down to rows 18-20, and update the Inference sheet cells to
# match (which ripples automatically through Conv/ReLU/Flatten/Dense/
# Softmax/Classification via recalculation).

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInference = $wb.Worksheets.Item("Inference")

# --- Update the Inference sheet values + highlight fill -------------------
# New inference window: B2=0 C2=1 D2=1 / B3=1 C3=1 D3=0
# (B2/B4/C4/D4/A* unchanged)
$wsInference.Range("C2").Value = 1
$wsInference.Range("D2").Value = 1
$wsInference.Range("B3").Value = 1
$wsInference.Range("C3").Value = 1
$wsInference.Range("D3").Value = 0

# Match the green highlight fill used by the corresponding filter window on
# the Inputs sheet (rows 18-20) by copying its cell formatting over, so the
# same shared style entries get reused instead of new ones being minted.
$wsInputs.Range("C18").Copy() | Out-Null
$wsInference.Range("C2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsInputs.Range("D18").Copy() | Out-Null
$wsInference.Range("D2").PasteSpecial(-4122) | Out-Null

$wsInputs.Range("B19").Copy() | Out-Null
$wsInference.Range("B3").PasteSpecial(-4122) | Out-Null

$wsInputs.Range("C19").Copy() | Out-Null
$wsInference.Range("C3").PasteSpecial(-4122) | Out-Null

# D3 no longer gets the highlight fill - clear it back to the default
# (unstyled) look by copying the format from an already-unstyled cell.
$wsInference.Range("B2").Copy() | Out-Null
$wsInference.Range("D3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View-state: move the selection to the new highlighted windows --------
$wsInputs.Activate() | Out-Null
$wsInputs.Range("B18:D20").Select() | Out-Null

$wsInference.Activate() | Out-Null
$wsInference.Range("D3").Select() | Out-Null
